$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.678.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.116.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  +1.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '

$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5259'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.10'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09093'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.171'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.119.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.818'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("E15").Value = '  +3.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '97.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.10%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001164'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.29%  '

$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06711'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '

$ws.Range("E21").Value = '  +1.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.413'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.764.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.69%  '

$ws.Range("E24").Value = '  +3.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.375'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.365.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.543'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.21%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.382'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.24%  '

$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.933'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02656'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06851'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.40%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2324'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6872'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.81%  '

$ws.Range("E43").Value = '  +0.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6437'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.313'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000369'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +14.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.706'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.256'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07322'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '82.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '
